$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing 2014 count value
$ws.Range("B5").Value = 432

# Add new row for 2015.
# A plain Range.Value assignment of a numeric-looking string like "2015"
# gets auto-coerced to a number by Excel's input parser (and forcing it to
# stay text via a leading apostrophe / "@" text format would introduce a new
# quotePrefix cell style that the original file never had). Instead, compute
# the text via a formula (so it is unambiguously a string) on a scratch cell,
# then Copy / PasteSpecial-values it into place: pasting values carries over
# the already-resolved string type without re-running the "looks like a
# number" autodetection, so the cell ends up as a plain shared-string cell
# with no extra styling - matching how the original rows were authored.
$scratch = $ws.Range("ZZ1")
$scratch.Formula = "=""2015"""
$scratch.Copy()
$ws.Range("A6").PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = 0

$ws.Range("B6").Value = 292
